$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.477.69'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '2.045.08'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.27'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.99'
$ws.Range("E8").Value = '  -7.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '62.76'
$ws.Range("E9").Value = '  +5.28%  '
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("E11").Value = '  -5.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.106'
$ws.Range("E12").Value = '  -3.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.941'
$ws.Range("E13").Value = '  +6.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.56'
$ws.Range("E14").Value = '  -5.15%  '
$ws.Range("D15").Value = '2.345.59'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("E16").Value = '  -5.19%  '
$ws.Range("D17").Value = '2.059.50'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").Value = '36.418.36'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.94'
$ws.Range("E19").Value = '  -6.09%  '
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").Value = '0.0₃0852'
$ws.Range("E21").Value = '  -5.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.31'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.17'
$ws.Range("E23").Value = '  -4.97%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.58'
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.12'
$ws.Range("E28").Value = '  -13.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.79'
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.01'
$ws.Range("E31").Value = '  -10.06%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  +4.74%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0589'
$ws.Range("E33").Value = '  -4.93%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.39'
$ws.Range("E34").Value = '  -7.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0877'
$ws.Range("E35").Value = '  +6.78%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("E38").Value = '  -6.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.01'
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.22'
$ws.Range("E40").Value = '  -8.07%  '
$ws.Range("E41").Value = '  -4.87%  '
$ws.Range("E42").Value = '  -5.12%  '
$ws.Range("E43").Value = '  -5.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.37'
$ws.Range("E44").Value = '  -4.24%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0898'
$ws.Range("E45").Value = '  -6.47%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.392.86'
$ws.Range("E46").Value = '  +6.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.70'
$ws.Range("E47").Value = '  -7.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.40'
$ws.Range("E48").Value = '  +9.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.91'
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("E50").Value = '  -5.32%  '
$ws.Range("D51").Value = '2.227.52'
$ws.Range("E51").Value = '  -0.47%  '
